$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the "Temps consacre" column to make room for
# a "Date" column; this shifts E->F, F->G, G->H, H->I, I->J, J->K.
$ws.Range("E1").EntireColumn.Insert()

# The old filler rows (9-19) are no longer needed.
$ws.Range("9:19").EntireRow.Delete()

# Clean up the blank filler cells that the column insert carried over into
# row 6 (former G6/J6 spacer cells, now H6/K6).
$ws.Range("H6").Clear()
$ws.Range("K6").Clear()

# New "Date" header cell.
$ws.Range("E6").Value = "Date"

# First existing task entry gets a date (1/20/2023).
$ws.Range("E7").NumberFormat = "mm-dd-yy"
$ws.Range("E7").Value = "1/20/2023"

# New second task entry: Kevin Carufel worked 8 hours on 1/22/2023.
$ws.Range("C8").Clear()
$ws.Range("C8").Value = "Kevin Carufel"

$ws.Range("E8").NumberFormat = "mm-dd-yy"
$ws.Range("E8").Value = "1/22/2023"

$ws.Range("F8").Value = 8

# Remove the now-unused filler cells around the new row 8 content.
$ws.Range("B8").Clear()
$ws.Range("D8").Clear()
$ws.Range("G8").Clear()
